# ErgoLux Hindi (hi-IN) translation workbook update
#
# A new localized string triple is inserted into the "Tabla13" table
# (sheet "hi-IN") right after the "strChkPower" row (old sheet row 31),
# i.e. it becomes the new row 32:
#   File    = localization\strings
#   Key     = strWindowPos
#   Comment = In "settings" form, tab "User interface"
#   English = Remember window position and size on startup
#
# This pushes every following table row down by one (old row 32 "strDlgReset"
# becomes row 33, ..., old row 332 becomes row 333), which Rows.Insert()
# handles for us (values, styles, row heights, and the table's own row
# collection all shift automatically).
#
# The existing "strChkDlgPath" row (row 25) also gets its previously-blank
# "Comment" column filled in with the same new comment text
# (In "settings" form, tab "User interface"), re-using the shared string.
#
# Column D (Comment) is additionally widened slightly to fit the new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new table row at sheet row 32 (shifts 32..332 -> 33..333) ---
$ws.Rows.Item(32).Insert()

$ws.Cells.Item(32, 2).Value = "localization\strings"
$ws.Cells.Item(32, 3).Value = "strWindowPos"
$ws.Cells.Item(32, 4).Value = 'In "settings" form, tab "User interface"'
$ws.Cells.Item(32, 5).Value = "Remember window position and size on startup"

# --- Fill in the previously-empty Comment cell for the strChkDlgPath row ---
$ws.Cells.Item(25, 4).Value = 'In "settings" form, tab "User interface"'

# --- Grow the "Tabla13" table / its AutoFilter to cover the new row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:F204"))

# --- Widen column D (Comment) slightly to fit the new text ---
$ws.Columns.Item(4).ColumnWidth = 35.022135416666664
